$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = 3427
$ws.Range("C3").Value  = 3114
$ws.Range("C4").Value  = 2127
$ws.Range("C5").Value  = 1839
$ws.Range("C6").Value  = 1347
$ws.Range("C7").Value  = 700
$ws.Range("C8").Value  = 603
$ws.Range("C9").Value  = 494
$ws.Range("C10").Value = 493
$ws.Range("C11").Value = 464
